$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: key + label
$ws.Range("A1").Value = "kodepeminatan"
$ws.Range("A2").Value = "Kode Peminatan"

# Column B: key + label
$ws.Range("B1").Value = "namapeminatan"
$ws.Range("B2").Value = "Nama Peminatan"

# Column C: key + label
$ws.Range("C1").Value = "kelompokkeahlian"
$ws.Range("C2").Value = "Kelompok Keahlian"

# Column D: left blank (not used)
$ws.Range("D1").Value = ""
$ws.Range("D2").Value = ""

# Column E: key + label
$ws.Range("E1").Value = "kuota"
$ws.Range("E2").Value = "Kuota"

$ws.Range("D2").Select()
